# chore: adapt column header formatting to respective input file names (#7)
#
# The sheet "AHB-Diff" holds a side-by-side comparison between two AHB
# ("Anwendungshandbuch") format versions. Its header row used the generic
# suffixes "_old"/"_new"; this rewires them to carry the concrete format
# version instead ("_FV2310" for the left/old block, "_FV2404" for the
# right/new block), freezes the header row, and wraps the sheet's used
# range in a real Excel Table ("Table1") so the new headers double as the
# table's column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells -------------------------------------------
# Columns A-J carried the "<Label>_old" headers (the FV2310 / "old" format
# version), columns L-U carried "<Label>_new" (the FV2404 / "new" format
# version); column K is the standalone "diff" column and stays untouched.
$fv2310Cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$fv2404Cols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
$labels = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Range($fv2310Cols[$i] + "1").Value = $labels[$i] + "_FV2310"
    $ws.Range($fv2404Cols[$i] + "1").Value = $labels[$i] + "_FV2404"
}

# --- 2. Freeze the header row -----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table ----------------------------
# xlSrcRange = 1, header row present = xlYes (1) so the table column names
# are taken straight from the header cells just renamed above.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U59"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
